$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the existing date-cell style (A57) down to the new date cells so we
# reuse the workbook's existing "short date" style instead of creating a
# brand-new numFmt/cellXf entry.
$ws.Range("A57").Copy()
$ws.Range("A58:A59").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Add two new time-log entries (rows 58 and 59)
$ws.Range("A58").Value = 45250
$ws.Range("B58").Value = 4
$ws.Range("C58").Value = "github messed up with the merging with the main and master class. So my code didn’t run so I had to fix the issue. I did a lot of research and debugging, I had to delete the .idea folder, and it worked"

$ws.Range("A59").Value = 45252
$ws.Range("B59").Value = 3
$ws.Range("C59").Value = "the budget panel code got taken out because of the github issues, so I had to re-do it, but the code has issues in it"

# Update selection / active cell to mirror the workbook view change
$ws.Range("C59").Select()
